$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 1.25
$ws.Range("C1").Value = 1.25
$ws.Range("D1").Value = 1.25
$ws.Range("A2").Value = 0.05023296795605869
$ws.Range("B2").Value = 1.4860880493934758
$ws.Range("C2").Value = 1.4300784251427108
$ws.Range("D2").Value = 1.271247924445413
$ws.Range("E2").Value = -0.000010000000000004976
$ws.Range("A3").Value = 0.13615592072401358
$ws.Range("B3").Value = 1.6823611648084293
$ws.Range("C3").Value = 1.551387807736352
$ws.Range("D3").Value = 1.3082481160780843
$ws.Range("E3").Value = -0.000009999999999998153
$ws.Range("A4").Value = 0.24774535458654837
$ws.Range("B4").Value = 1.86647929190935
$ws.Range("C4").Value = 1.6600384839000062
$ws.Range("D4").Value = 1.3584889219028982
$ws.Range("A5").Value = 0.3395483796244561
$ws.Range("B5").Value = 2.018148261365888
$ws.Range("C5").Value = 1.748130336594594
$ws.Range("D5").Value = 1.4178582624545673
$ws.Range("E5").Value = 0.21972549121774831
$ws.Range("A6").Value = 0.37489509211216554
$ws.Range("B6").Value = 2.1110528095156904
$ws.Range("C6").Value = 1.810713896064194
$ws.Range("D6").Value = 1.4802457716900683
$ws.Range("E6").Value = 0.5346417644211807
$ws.Range("A7").Value = 0.3596941945957446
$ws.Range("B7").Value = 2.144898199484408
$ws.Range("C7").Value = 1.8485298505934096
$ws.Range("D7").Value = 1.540064490191559
$ws.Range("E7").Value = 0.7224931491655694
$ws.Range("A8").Value = 0.3154747180195981
$ws.Range("B8").Value = 2.1372998490718156
$ws.Range("C8").Value = 1.8681783761752564
$ws.Range("D8").Value = 1.593815672490782
$ws.Range("E8").Value = 0.793762503223617
$ws.Range("A9").Value = 0.2632708888319195
$ws.Range("B9").Value = 2.110025665677201
$ws.Range("C9").Value = 1.8778486451314116
$ws.Range("D9").Value = 1.6402745718775598
$ws.Range("E9").Value = 0.7897200391716307
$ws.Range("A10").Value = 0.21668730322318777
$ws.Range("B10").Value = 2.08003499098652
$ws.Range("C10").Value = 1.8841029198827006
$ws.Range("D10").Value = 1.6798530095973663
$ws.Range("E10").Value = 0.7537867837257558
$ws.Range("A11").Value = 0.18119790590537352
$ws.Range("B11").Value = 2.0565187109776675
$ws.Range("C11").Value = 1.890698760611018
$ws.Range("D11").Value = 1.7137529227215933
$ws.Range("E11").Value = 0.7166163498849444
$ws.Range("A12").Value = 0.15651849202618223
$ws.Range("B12").Value = 2.0419621423831495
$ws.Range("C12").Value = 1.8988942557461408
$ws.Range("D12").Value = 1.7432917524911336
$ws.Range("E12").Value = 0.6929377683064195
$ws.Range("A13").Value = 0.13963351566524013
$ws.Range("B13").Value = 2.0348248669860265
$ws.Range("C13").Value = 1.9083927339788556
$ws.Range("D13").Value = 1.7695297327956738
$ws.Range("E13").Value = 0.6851599791081456
$ws.Range("A14").Value = 0.1271384339293077
$ws.Range("B14").Value = 2.032128302320473
$ws.Range("C14").Value = 1.918292513151762
$ws.Range("D14").Value = 1.7931636040529058
$ws.Range("E14").Value = 0.6890706978036171
$ws.Range("A15").Value = 0.11649034860870687
$ws.Range("B15").Value = 2.0311588089722292
$ws.Range("C15").Value = 1.9277319863107136
$ws.Range("D15").Value = 1.8145831724956447
$ws.Range("E15").Value = 0.6986951553308774
$ws.Range("A16").Value = 0.10633028113800491
$ws.Range("B16").Value = 2.0301952427508
$ws.Range("C16").Value = 1.9361803478724735
$ws.Range("D16").Value = 1.833988258818609
$ws.Range("E16").Value = 0.7091798039670456
$ws.Range("A17").Value = 0.09625747084677917
$ws.Range("B17").Value = 2.0285347448974242
$ws.Range("C17").Value = 1.943465296423371
$ws.Range("D17").Value = 1.8514974425657023
$ws.Range("E17").Value = 0.71779917610711
$ws.Range("A18").Value = 0.08641533432093404
$ws.Range("B18").Value = 2.0261755097548235
$ws.Range("C18").Value = 1.949666718797797
$ws.Range("D18").Value = 1.8672184686127231
$ws.Range("E18").Value = 0.723741846093582
$ws.Range("A19").Value = 0.07712176853005281
$ws.Range("B19").Value = 2.023439053327615
$ws.Range("C19").Value = 1.9549813791549262
$ws.Range("D19").Value = 1.8812783212370634
$ws.Range("E19").Value = 0.7273846139071285
$ws.Range("A20").Value = 0.06864030360422026
$ws.Range("B20").Value = 2.020686306237904
$ws.Range("C20").Value = 1.9596180506185361
$ws.Range("D20").Value = 1.893825039887139
$ws.Range("E20").Value = 0.7295590335177972
$ws.Range("A21").Value = 0.06109234061456484
$ws.Range("B21").Value = 2.0181705650975434
$ws.Range("C21").Value = 1.9637414923014205
$ws.Range("D21").Value = 1.9050161371560754
$ws.Range("E21").Value = 0.7310528983591983
$ws.Range("A22").Value = 0.054464762812490124
$ws.Range("B22").Value = 2.016003341115631
$ws.Range("C22").Value = 1.9674581151009236
$ws.Range("D22").Value = 1.9150049855124354
$ws.Range("E22").Value = 0.7323858649004814
$ws.Range("A23").Value = 0.04866017069064739
$ws.Range("B23").Value = 2.0141849576929975
$ws.Range("C23").Value = 1.970826375015069
$ws.Range("D23").Value = 1.923931183008686
$ws.Range("E23").Value = 0.7337887084134538
$ws.Range("A24").Value = 0.04355079287463968
$ws.Range("B24").Value = 2.0126559726964794
$ws.Range("C24").Value = 1.973875690377376
$ws.Range("D24").Value = 1.9319164140805873
$ws.Range("E24").Value = 0.7352864325134242
$ws.Range("A25").Value = 0.03901634314508759
$ws.Range("B25").Value = 2.011341605976377
$ws.Range("C25").Value = 1.976623466379496
$ws.Range("D25").Value = 1.9390646813512085
$ws.Range("E25").Value = 0.7368038297165093
$ws.Range("A26").Value = 0.03496169169632725
$ws.Range("B26").Value = 2.010178403291931
$ws.Range("C26").Value = 1.9790856554771294
$ws.Range("D26").Value = 1.9454649163258735
$ws.Range("E26").Value = 0.7382475202232286
$ws.Range("A27").Value = 0.03131893366158255
$ws.Range("B27").Value = 2.00912334985576
$ws.Range("C27").Value = 1.981281098680161
$ws.Range("D27").Value = 1.9511941753435633
$ws.Range("E27").Value = 0.7395499137402098
$ws.Range("A28").Value = 0.028041202681653634
$ws.Range("B28").Value = 2.008151973491249
$ws.Range("C28").Value = 1.983231075340254
$ws.Range("D28").Value = 1.956320377176855
$ws.Range("E28").Value = 0.7406809699856729
$ws.Range("A29").Value = 0.02509430702464043
$ws.Range("B29").Value = 2.0072506551031344
$ws.Range("C29").Value = 1.9849574804844021
$ws.Range("D29").Value = 1.96090410219022
$ws.Range("E29").Value = 0.7416400516759384
$ws.Range("A30").Value = 0.022450047579245548
$ws.Range("B30").Value = 2.006410086325446
$ws.Range("C30").Value = 1.986479320850795
$ws.Range("D30").Value = 1.9649996407623904
$ws.Range("E30").Value = 0.742441770323128
$ws.Range("A31").Value = 0.020082084534033404
$ws.Range("B31").Value = 2.0056177533874293
$ws.Range("C31").Value = 1.9878121110131841
$ws.Range("D31").Value = 1.968655270898644
$ws.Range("E31").Value = 0.7431011325011929
$ws.Range("A32").Value = 0.017964774108568473
$ws.Range("B32").Value = 2.0048577522061723
$ws.Range("C32").Value = 1.9889632743296846
$ws.Range("D32").Value = 1.9719134942163214
$ws.Range("E32").Value = 0.7436263274707845
$ws.Range("A33").Value = 0.01607269506893391
$ws.Range("B33").Value = 2.0041021641944976
$ws.Range("C33").Value = 1.9899361958772046
$ws.Range("D33").Value = 1.9748104745143573
$ws.Range("E33").Value = 0.7440095381020904
$ws.Range("A34").Value = 0.01438241251262249
$ws.Range("B34").Value = 2.003320608171677
$ws.Range("C34").Value = 1.9907181909738854
$ws.Range("D34").Value = 1.9773763865435163
$ws.Range("E34").Value = 0.7442314617025055
$ws.Range("A35").Value = 0.012871206927983917
$ws.Range("B35").Value = 2.0024550435599924
$ws.Range("C35").Value = 1.9912967916874118
$ws.Range("D35").Value = 1.9796334656749992
$ws.Range("E35").Value = 0.7442445248055773
$ws.Range("A36").Value = 0.011521105315361876
$ws.Range("B36").Value = 2.0014550705763408
$ws.Range("C36").Value = 1.991619969778408
$ws.Range("D36").Value = 1.98159741011612
$ws.Range("E36").Value = 0.7439935887625897
$ws.Range("A37").Value = 0.010312147221051145
$ws.Range("B37").Value = 2.0001935315109174
$ws.Range("C37").Value = 1.9916553450668706
$ws.Range("D37").Value = 1.9832710610416517
$ws.Range("E37").Value = 0.7433591921184995
$ws.Range("A38").Value = 0.009234843140444626
$ws.Range("B38").Value = 1.998594495066222
$ws.Range("C38").Value = 1.9912574741861784
$ws.Range("D38").Value = 1.9846501701038632
$ws.Range("E38").Value = 0.7422323557595297
$ws.Range("A39").Value = 0.008266426136566208
$ws.Range("B39").Value = 1.996351595433337
$ws.Range("C39").Value = 1.9903811422307083
$ws.Range("D39").Value = 1.9857032983835157
$ws.Range("E39").Value = 0.7403269093112094
$ws.Range("A40").Value = 0.007414340215884392
$ws.Range("B40").Value = 1.9933880413703489
$ws.Range("C40").Value = 1.9886349323570018
$ws.Range("D40").Value = 1.9863949252523307
$ws.Range("E40").Value = 0.7374557220580851
$ws.Range("A41").Value = 0.006634767995785583
$ws.Range("B41").Value = 1.9889071145083963
$ws.Range("C41").Value = 1.9860391540074187
$ws.Range("D41").Value = 1.9866210222853766
$ws.Range("E41").Value = 0.7329093014085913
$ws.Range("A42").Value = 0.005983954064531185
$ws.Range("B42").Value = 1.983022203373304
$ws.Range("C42").Value = 1.9815092052260121
$ws.Range("D42").Value = 1.9862971285832902
$ws.Range("E42").Value = 0.7264549160481463
$ws.Range("A43").Value = 0.005336748184002808
$ws.Range("B43").Value = 1.973524429505503
$ws.Range("C43").Value = 1.9755217153616864
$ws.Range("D43").Value = 1.9851475856662892
$ws.Range("E43").Value = 0.7162483076096527
$ws.Range("A44").Value = 0.004910787922650093
$ws.Range("B44").Value = 1.961605233168871
$ws.Range("C44").Value = 1.964926039821358
$ws.Range("D44").Value = 1.9830287739415215
$ws.Range("E44").Value = 0.7023841823938852
$ws.Range("A45").Value = 0.004292172849618721
$ws.Range("B45").Value = 1.9407824658083943
$ws.Range("C45").Value = 1.952264664906494
$ws.Range("D45").Value = 1.97922660620954
$ws.Range("E45").Value = 0.6798003660716734
$ws.Range("A46").Value = 0.004249885870871826
$ws.Range("B46").Value = 1.9168809277971521
$ws.Range("C46").Value = 1.9279438068554464
$ws.Range("D46").Value = 1.9736154951524254
$ws.Range("E46").Value = 0.6507640613044322
$ws.Range("A47").Value = 0.0033622872044484758
$ws.Range("B47").Value = 1.869944065287759
$ws.Range("C47").Value = 1.9026933729113937
$ws.Range("D47").Value = 1.9642850664646052
$ws.Range("E47").Value = 0.6005918748043716
$ws.Range("A48").Value = 0.00429753512588564
$ws.Range("B48").Value = 1.8237669065864215
$ws.Range("C48").Value = 1.8461276765741315
$ws.Range("D48").Value = 1.9516384320755686
$ws.Range("E48").Value = 0.5412923681002331
$ws.Range("A49").Value = 0.002133592250105823
$ws.Range("B49").Value = 1.7145702052066294
$ws.Range("C49").Value = 1.799532792285102
$ws.Range("D49").Value = 1.9303022916573642
$ws.Range("E49").Value = 0.42812167272156215
$ws.Range("A50").Value = 0.006098161827261275
$ws.Range("B50").Value = 1.632490986776422
$ws.Range("C50").Value = 1.663779206528666
$ws.Range("D50").Value = 1.9034992742180794
$ws.Range("E50").Value = 0.3115603295652709
$ws.Range("A51").Value = -0.0008690252248051982
$ws.Range("B51").Value = 1.3680855517747998
$ws.Range("C51").Value = 1.5906655833900007
$ws.Range("D51").Value = 1.8553120391981843
$ws.Range("E51").Value = 0.04974467181940624
$ws.Range("B52").Value = 1.25
$ws.Range("C52").Value = 1.25
$ws.Range("D52").Value = 1.25
